$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current row 53, pushing the existing
# rows 53-54 down to 55-56 (their contents remain unchanged).
$ws.Rows.Item(53).Insert()
$ws.Rows.Item(53).Insert()

# New row 53: weekly update for "Americana (o)" variety from Provincia de Limari.
$ws.Range("A53").Value = 11
$ws.Range("B53").Value = "Vega Monumental Concepción"
$ws.Range("C53").Value = "Bíobío"
$ws.Range("D53").Value = 44509
$ws.Range("E53").Value = 8
$ws.Range("F53").Value = 100112021
$ws.Range("G53").Value = "Ají"
$ws.Range("H53").Value = "Americana (o)"
$ws.Range("I53").Value = "Primera"
$ws.Range("J53").Value = 50
$ws.Range("K53").Value = 32000
$ws.Range("L53").Value = 34000
$ws.Range("M53").Value = 32800
$ws.Range("N53").Value = "$/caja 25 kilos"
$ws.Range("O53").Value = "Provincia de Limarí"
$ws.Range("P53").Value = 1312
$ws.Range("Q53").Value = 25
$ws.Range("R53").Value = "Hortaliza"

# New row 54: weekly update for the existing "Inferno" variety, same week as row 53.
$ws.Range("A54").Value = 11
$ws.Range("B54").Value = "Vega Monumental Concepción"
$ws.Range("C54").Value = "Bíobío"
$ws.Range("D54").Value = 44509
$ws.Range("E54").Value = 8
$ws.Range("F54").Value = 100112021
$ws.Range("G54").Value = "Ají"
$ws.Range("H54").Value = "Inferno"
$ws.Range("I54").Value = "Primera"
$ws.Range("J54").Value = 40
$ws.Range("K54").Value = 23000
$ws.Range("L54").Value = 24000
$ws.Range("M54").Value = 23500
$ws.Range("N54").Value = "$/caja 12 kilos"
$ws.Range("O54").Value = "Región de Arica y Parinacota"
$ws.Range("P54").Value = 1958
$ws.Range("Q54").Value = 12
$ws.Range("R54").Value = "Hortaliza"

# Ensure the date cells keep the date/time number format used elsewhere
# in column D (style index 2 in the original workbook).
$ws.Range("D53").NumberFormat = $ws.Range("D52").NumberFormat
$ws.Range("D54").NumberFormat = $ws.Range("D52").NumberFormat
